$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The shop-order of the 3 "light dragon" skin rows (76-78) is corrected
# (shopOrder 1,2,3 was sitting on the wrong rows) and the temporary
# red/bold "not localized yet" highlight that sat on row 78's [powerup]
# cell is removed (the commit message: "TEXTS not localized yet").
# Values below are the final, authoritative cell contents lifted
# straight from the target worksheet.
# -----------------------------------------------------------------------

# Row 76 ("transform_gold_LOW" / shopOrder 3)
$ws.Range("E76").Value = "transform_gold_LOW"
$ws.Range("F76").Value = 3
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 70
$ws.Range("I76").Value = 7
$ws.Range("Q76").Value = "TID_SKIN_LIGHT_3_NAME"
$ws.Range("R76").Value = "TID_DRAGON_LIGHT_3_DESC"

# Row 77 ("disguise_hp" / shopOrder 1)
$ws.Range("E77").Value = "disguise_hp"
$ws.Range("F77").Value = 1
$ws.Range("G77").Value = 900000
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 1
$ws.Range("Q77").Value = "TID_SKIN_LIGHT_1_NAME"
$ws.Range("R77").Value = "TID_DRAGON_LIGHT_1_DESC"

# Row 78 ("disguise_furyDuration_hp" / shopOrder 2)
$ws.Range("E78").Value = "disguise_furyDuration_hp"
$ws.Range("F78").Value = 2
$ws.Range("H78").Value = 65
$ws.Range("I78").Value = 4
$ws.Range("Q78").Value = "TID_SKIN_LIGHT_2_NAME"
$ws.Range("R78").Value = "TID_DRAGON_LIGHT_2_DESC"

# Row 78's [powerup] cell (E78) loses the special red/bold "TODO"
# formatting - pull plain formatting from E77 (same column, same visual
# style as the rest of the table) so the leftover highlight style is no
# longer referenced by any cell.
$ws.Range("E77").Copy()
$ws.Range("E78").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column I ([unlockLevel]/trailing column) widens slightly now that it
# stands on its own (was merged in a min=8 max=9 column-width span).
$ws.Columns.Item(9).ColumnWidth = 11.5

# View state: scrolled up one row, selection moved from column O to I.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I78").Select()
